# Bugs.xlsx edit: add two new bug rows (sumby/sumfig NA check, and
# .reserved warning in clean_up()), format their Description cells with
# the "Bad" (red) cell style, and update the tracked selection.
#
# Commit message: "checking of code-tree existing in source was failling
# always - code error now fixed"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook window was minimized in the source edit.
$excel.WindowState = -4140   # xlMinimized

# ---- Row 46: new bug raised by SB on 2019-10-16 (serial 43754) ----
# Copy the date format from the existing date cell above (B45) so the new
# date cell reuses the same style instead of creating a duplicate one.
$ws.Cells.Item(45, 2).Copy()
$ws.Cells.Item(46, 2).PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(46, 1).Value = "SB"
$ws.Cells.Item(46, 2).Value = 43754
$ws.Cells.Item(46, 3).Value = "If sumby is fed an empty or all NA variable, then the sumfig crashes. Add a check in and disable the figure."
$ws.Cells.Item(46, 3).Style = "Bad"

# ---- Row 47: new bug raised by SB on 2019-10-16 (serial 43754) ----
$ws.Cells.Item(45, 2).Copy()
$ws.Cells.Item(47, 2).PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(47, 1).Value = "SB"
$ws.Cells.Item(47, 2).Value = 43754
$ws.Cells.Item(47, 3).Value = "Give a better warning if .reserved is undefined in clean_up()"
$ws.Cells.Item(47, 3).Style = "Bad"

# Update the tracked selection/active cell to C36 (as in the source edit).
$ws.Range("C36").Select()
